$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 9.014599999999996
$ws.Range("B6").Value = 6.663799999999998
$ws.Range("B7").Value = 5.074900000000003
$ws.Range("D7").Value = -7.606599999999991
$ws.Range("B8").Value = 6.7071
$ws.Range("D11").Value = -7.915699999999999
$ws.Range("D12").Value = -7.239399999999996
$ws.Range("D15").Value = -8.578299999999993
$ws.Range("B16").Value = 7.070299999999996
$ws.Range("B20").Value = 9.564699999999991
$ws.Range("D20").Value = -7.934899999999996
$ws.Range("B21").Value = 9.347099999999998
$ws.Range("D21").Value = -8.167699999999996
$ws.Range("D22").Value = -7.4768
$ws.Range("D23").Value = -7.250699999999997
$ws.Range("B28").Value = 6.239799999999998
$ws.Range("B29").Value = 5.214900000000003
$ws.Range("D29").Value = -7.309099999999996
$ws.Range("B30").Value = 5.48
$ws.Range("B32").Value = 7.373999999999996
$ws.Range("D34").Value = -7.659899999999999
$ws.Range("B40").Value = 9.091899999999995
$ws.Range("D42").Value = -8.037400000000002
$ws.Range("D43").Value = -8.324999999999996
$ws.Range("D44").Value = -7.756900000000001
$ws.Range("D45").Value = -7.930499999999998
$ws.Range("B46").Value = 6.1656
$ws.Range("D46").Value = -8.162199999999999
$ws.Range("D50").Value = -8.202500000000001
$ws.Range("B51").Value = 5.3279
$ws.Range("D51").Value = -7.706099999999991
$ws.Range("B52").Value = 5.377999999999999
$ws.Range("B57").Value = 5.626799999999996
$ws.Range("D57").Value = -7.903099999999998
$ws.Range("B59").Value = 5.032199999999999
$ws.Range("B62").Value = 5.797399999999999
$ws.Range("D65").Value = -7.736199999999997
$ws.Range("B66").Value = 5.382600000000001
$ws.Range("D66").Value = -7.283
$ws.Range("D67").Value = -6.439300000000004
$ws.Range("B73").Value = 8.403900000000002
$ws.Range("B74").Value = 9.334399999999993
$ws.Range("B77").Value = 8.887800000000009
$ws.Range("D79").Value = -6.363900000000005
$ws.Range("D84").Value = -9.105200000000004
$ws.Range("D87").Value = -8.144899999999998
$ws.Range("B92").Value = 4.845899999999999
$ws.Range("D92").Value = -6.420999999999999
$ws.Range("D97").Value = -8.662900000000004
$ws.Range("B100").Value = 5.721999999999993
